$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (D) and volume-change (E) columns; B/C updated for the
# 3-way reshuffle of the Toncoin / InjectiveProtocol / Kaspa rows (29-31).
# A leading apostrophe is used for D-column values that would otherwise be
# auto-parsed as numbers by Excel, so they stay stored as text (matching the
# source data, which keeps these as plain strings).

$ws.Range('D2').Value = '47.331.70'
$ws.Range('D3').Value = '2.492.45'
$ws.Range('E3').Value = '  +2.78%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '''323.00'
$ws.Range('E5').Value = '  +1.35%  '
$ws.Range('D6').Value = '''108.19'
$ws.Range('E6').Value = '  +5.54%  '
$ws.Range('D7').Value = '''0.526'
$ws.Range('E7').Value = '  +2.31%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '''0.544'
$ws.Range('E9').Value = '  +3.35%  '
$ws.Range('D10').Value = '''38.11'
$ws.Range('E10').Value = '  +7.32%  '
$ws.Range('D11').Value = '''0.0813'
$ws.Range('E11').Value = '  +1.72%  '
$ws.Range('E12').Value = '  +1.44%  '
$ws.Range('D13').Value = '''18.42'
$ws.Range('E13').Value = '  +0.99%  '
$ws.Range('D14').Value = '''7.21'
$ws.Range('E14').Value = '  +1.98%  '
$ws.Range('D15').Value = '2.886.12'
$ws.Range('E15').Value = '  +2.88%  '
$ws.Range('D16').Value = '2.498.42'
$ws.Range('E16').Value = '  +4.81%  '
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('D18').Value = '47.264.25'
$ws.Range('E18').Value = '  +4.65%  '
$ws.Range('D19').Value = '''12.90'
$ws.Range('E19').Value = '  +5.69%  '
$ws.Range('E20').Value = '  +5.41%  '
$ws.Range('D21').Value = '0.0₃0941'
$ws.Range('E21').Value = '  +2.14%  '
$ws.Range('D22').Value = '''70.67'
$ws.Range('E22').Value = '  +2.53%  '
$ws.Range('D23').Value = '''2.44'
$ws.Range('E23').Value = '  +7.14%  '
$ws.Range('D24').Value = '''250.78'
$ws.Range('E24').Value = '  +2.45%  '
$ws.Range('E25').Value = '  +4.41%  '
$ws.Range('D26').Value = '''26.19'
$ws.Range('E26').Value = '  +1.67%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.21'
$ws.Range('E29').Value = '  -3.28%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '''34.99'
$ws.Range('E30').Value = '  +6.32%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').Value = '''0.138'
$ws.Range('E31').Value = '  +10.41%  '
$ws.Range('D32').Value = '''49.34'
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('D33').Value = '''5.49'
$ws.Range('E33').Value = '  +5.34%  '
$ws.Range('D34').Value = '''19.69'
$ws.Range('E34').Value = '  -1.75%  '
$ws.Range('D35').Value = '''0.0791'
$ws.Range('E35').Value = '  +3.75%  '
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('E37').Value = '  +6.22%  '
$ws.Range('D38').Value = '''4.68'
$ws.Range('E38').Value = '  +5.60%  '
$ws.Range('E39').Value = '  +3.86%  '
$ws.Range('E40').Value = '  +2.04%  '
$ws.Range('E41').Value = '  +2.10%  '
$ws.Range('D42').Value = '''121.76'
$ws.Range('E42').Value = '  -3.36%  '
$ws.Range('D43').Value = '''21.19'
$ws.Range('E43').Value = '  +2.22%  '
$ws.Range('E44').Value = '  +3.41%  '
$ws.Range('D45').Value = '1.965.25'
$ws.Range('E45').Value = '  +1.54%  '
$ws.Range('D46').Value = '''3.02'
$ws.Range('E46').Value = '  +2.50%  '
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('D49').Value = '''9.07'
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('E50').Value = '  +9.94%  '
$ws.Range('D51').Value = '''79.33'
$ws.Range('E51').Value = '  +3.75%  '
